$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
}

# Row 2 - Bitcoin
Set-TextValue "D2" "60.312.05"
Set-TextValue "E2" "  -2.70%  "

# Row 3 - Ethereum
Set-TextValue "D3" "3.299.96"
Set-TextValue "E3" "  -3.47%  "

# Row 4 - TetherUSD
Set-TextValue "E4" "  +0.04%  "

# Row 5 - BNB
Set-TextValue "D5" "557.18"

# Row 6 - Solana
Set-TextValue "D6" "141.12"
Set-TextValue "E6" "  -8.33%  "

# Row 7 - USDC
Set-TextValue "E7" "  -0.10%  "

# Row 8 - LidoStakedEther
Set-TextValue "D8" "3.301.82"
Set-TextValue "E8" "  -3.42%  "

# Row 9 - XRP
Set-TextValue "E9" "  -3.59%  "

# Row 10 - Toncoin
Set-TextValue "E10" "  -1.31%  "

# Row 11 - Dogecoin
Set-TextValue "E11" "  -4.98%  "

# Row 12 - Cardano
Set-TextValue "E12" "  -2.65%  "

# Row 13 - WrappedliquidstakedEther2.0
Set-TextValue "D13" "3.865.07"
Set-TextValue "E13" "  -3.46%  "

# Row 14 - TRON
Set-TextValue "E14" "  -0.23%  "

# Row 15 - Avalanche
Set-TextValue "D15" "26.68"
Set-TextValue "E15" "  -6.83%  "

# Row 16 - WrappedEther
Set-TextValue "D16" "3.301.10"
Set-TextValue "E16" "  -3.49%  "

# Row 17 - ShibaInu
Set-TextValue "D17" "0.0000164"
Set-TextValue "E17" "  -4.87%  "

# Row 18 - WrappedBTC
Set-TextValue "D18" "60.281.00"
Set-TextValue "E18" "  -2.78%  "

# Row 19 - Polkadot
Set-TextValue "D19" "6.06"
Set-TextValue "E19" "  -7.78%  "

# Row 20 - Chainlink
Set-TextValue "D20" "13.72"
Set-TextValue "E20" "  -4.88%  "

# Row 21 - Uniswap
Set-TextValue "D21" "8.51"
Set-TextValue "E21" "  -4.86%  "

# Row 22 - BitcoinCash
Set-TextValue "D22" "373.54"
Set-TextValue "E22" "  -2.26%  "

# Row 23 - Dai
Set-TextValue "E23" "  +0.08%  "

# Row 24 - Litecoin
Set-TextValue "D24" "72.47"
Set-TextValue "E24" "  -4.72%  "

# Row 25 - Polygon
Set-TextValue "E25" "  -6.70%  "

# Row 26 - WrappedeETH
Set-TextValue "D26" "3.431.48"
Set-TextValue "E26" "  -3.65%  "

# Row 27 - PEPE
Set-TextValue "D27" "0.0000103"
Set-TextValue "E27" "  -9.08%  "

# Row 28 - Kaspa
Set-TextValue "E28" "  -2.16%  "

# Row 29 - Binance-PegBSC-USD
Set-TextValue "E29" "  +0.43%  "

# Row 30 - RenderToken
Set-TextValue "D30" "7.05"
Set-TextValue "E30" "  -7.72%  "

# Row 31 - USDe
Set-TextValue "D31" "1.00"
Set-TextValue "E31" "  +0.05%  "

# Row 32 - PancakeSwap
Set-TextValue "E32" "  -4.99%  "

# Row 33 - InternetComputer(DFINITY)
Set-TextValue "D33" "7.42"
Set-TextValue "E33" "  -5.87%  "

# Row 34 - EthereumClassic
Set-TextValue "D34" "22.54"
Set-TextValue "E34" "  -3.08%  "

# Row 35 - Fetch.AI
Set-TextValue "D35" "1.24"
Set-TextValue "E35" "  -6.39%  "

# Row 36 and 37 swap: Monero <-> NEARProtocol (with updated figures)
Set-TextValue "B36" "NEARProtocol"
Set-TextValue "C36" "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue "D36" "5.04"
Set-TextValue "E36" "  -9.37%  "

Set-TextValue "B37" "Monero"
Set-TextValue "C37" "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D37" "166.10"
Set-TextValue "E37" "  -1.29%  "

# Row 38 - ImmutableX
Set-TextValue "E38" "  -4.67%  "

# Row 40 - RenzoRestakedETH
Set-TextValue "D40" "3.331.12"
Set-TextValue "E40" "  -3.57%  "

# Row 41 - Hedera
Set-TextValue "E41" "  -7.83%  "

# Row 42 - EnergySwap
Set-TextValue "D42" "25.56"
Set-TextValue "E42" "  -17.57%  "

# Row 43 - OKB
Set-TextValue "D43" "41.64"
Set-TextValue "E43" "  -2.56%  "

# Row 44 - Mantle
Set-TextValue "D44" "0.748"
Set-TextValue "E44" "  -4.14%  "

# Row 45 - ONDO
Set-TextValue "E45" "  -3.80%  "

# Row 46 - Filecoin
Set-TextValue "D46" "4.09"
Set-TextValue "E46" "  -7.51%  "

# Row 47 - Stacks
Set-TextValue "E47" "  -6.29%  "

# Row 48 - FirstDigitalUSD
Set-TextValue "E48" "  -0.01%  "

# Row 49 - Maker
Set-TextValue "D49" "2.320.99"
Set-TextValue "E49" "  -8.97%  "

# Row 50 - InjectiveProtocol
Set-TextValue "D50" "21.57"
Set-TextValue "E50" "  -6.76%  "

# Row 51 - Cosmos
Set-TextValue "D51" "6.33"
Set-TextValue "E51" "  -6.99%  "
